$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the target paragraph:
#   "Coeficiente de transferencia de calor por convección,  h = 4,8
#    w/(m2.oC)"
# and collapse the double space after the comma:
#   "convección,  h"  ->  "convección, h"
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Coeficiente de transferencia de calor por*") {
        $target = $cand
        break
    }
}

$pStart = $target.Range.Start
$text = $target.Range.Text
$idx = $text.IndexOf(",  h")
$rngStart = $pStart + $idx
$rngEnd = $rngStart + 4
$rng = $d.Range($rngStart, $rngEnd)
$rng.Text = ", h"

# ---------------------------------------------------------------------
# That text replacement merges adjacent runs that already shared the
# same formatting ("convección" + ", h", and separately
# " = " + "4,8" + " " + "w/(m"). Restore the original run boundaries by
# toggling Bold on/off across the tail of each run we need to peel off
# -- this forces a fresh run to start at the toggle point without
# changing any visible formatting, and merges cleanly into one run with
# whatever already-separate content lies beyond it (so a single toggle
# only ever introduces the one new boundary we want).
# ---------------------------------------------------------------------
function Split-Tail($fromPos, $toPos) {
    $r = $d.Range($fromPos, $toPos)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# Refresh paragraph/text after the edit (length shrank by 1 char).
$pStart = $target.Range.Start
$text = $target.Range.Text

# Split "convección" | ", h"
$splitPos = $pStart + $text.IndexOf("convección") + ("convección").Length
$endOfMergedRun = $pStart + $text.IndexOf(" = ")
Split-Tail $splitPos $endOfMergedRun

# Split " = 4,8 w/(m" -> " = " | "4,8" | " " | "w/(m"
# (processed right-to-left so each toggle only introduces one boundary)
$segStart = $pStart + $text.IndexOf(" = 4,8 w/(m")
$segEnd = $pStart + $text.IndexOf("2.")   # start of the superscript "2" run
$boundary1 = $segStart + 3                # after " = "
$boundary2 = $segStart + 3 + 3            # after "4,8"
$boundary3 = $segStart + 3 + 3 + 1        # after " "

Split-Tail $boundary3 $segEnd
Split-Tail $boundary2 $boundary3
Split-Tail $boundary1 $boundary2

# Split "C)" -> "C" | ")"
$pStart = $target.Range.Start
$text = $target.Range.Text
$cPos = $pStart + $text.IndexOf("C)") + 1
$paraEnd = $pStart + $text.Length
Split-Tail $cPos $paraEnd
